$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Noviembre")

# Row 3 - owner_crm
$ws.Range("B3").Value = 9524
$ws.Range("C3").Value = 734
$ws.Range("D3").Value = 248

# Row 6 - owner_krossboarder-remesas
$ws.Range("B6").Value = 6942
$ws.Range("E6").Value = 697

# Row 7 - owner_marketplace
$ws.Range("B7").Value = 16799
$ws.Range("F7").Value = 12529

# Row 8 - owner_promos
$ws.Range("B8").Value = 9682
$ws.Range("C8").Value = 627
$ws.Range("D8").Value = 358
$ws.Range("E8").Value = 1420

# Row 11 - retail
$ws.Range("B11").Value = 13257
$ws.Range("C11").Value = 849
$ws.Range("D11").Value = 343
$ws.Range("E11").Value = 1281
$ws.Range("F11").Value = 10783

# Row 12 - financial
$ws.Range("B12").Value = 15249
$ws.Range("E12").Value = 504

# Row 14 - product cx
$ws.Range("B14").Value = 9524
$ws.Range("C14").Value = 734
$ws.Range("D14").Value = 248
